$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# New header cells
$ws.Range("H7").Value = "AlarmLoadingDetail"
$ws.Range("I7").Value = "StandbyLoadingDetail"

# New data cells for rows 8-10
$ws.Range("H8").Value = "Battery Alarm (A)"
$ws.Range("I8").Value = "Battery Standby (A)"

$ws.Range("H9").Value = "Battery Alarm (A)"
$ws.Range("I9").Value = "Battery Standby (A)"

$ws.Range("H10").Value = "Battery Alarm (A)"
$ws.Range("I10").Value = "Battery Standby (A)"

$ws.Range("G7").Copy() | Out-Null
$ws.Range("H7:I7").PasteSpecial(-4122) | Out-Null

$ws.Range("B8").Copy() | Out-Null
$ws.Range("H8:I10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws.Columns.Item(8).ColumnWidth = 17.6640625
$ws.Columns.Item(9).ColumnWidth = 19.6640625
$ws.Columns.Item(8).BestFit = $true
$ws.Columns.Item(9).BestFit = $true

$ws.Range("I12").Select()
